$wb = $excel.ActiveWorkbook

# Sheet order in this workbook:
#   1 Funciones_Objetivo
#   2 Restricciones_del_lider
#   3 Restricciones_del_follower
#   4 Punto_modificado
#   5 Vector_bf
#   6 Vector_BF   (name differs from sheet 5 only by case -> use index, not name, to avoid
#                  a case-insensitive Worksheets.Item("...") collision)
#   7 Vector_Alpha

# --- Sheet: Restricciones_del_follower ---
$ws = $wb.Worksheets.Item(3)

$ws.Range("A2").Value = "4.49 - x - 0.5y"
$ws.Range("B2").Value = -2.49
$ws.Range("C2").Value = "J_0_L0_v"
$ws.Range("D2").Value = 0.62
$ws.Range("E2").Value = 0.4
$ws.Range("F2").Value = 0

$ws.Range("A3").Value = "-4.4125 - 0.25x + y"
$ws.Range("B3").Value = 2.4124999999999996
$ws.Range("C3").Value = "J_0_L0_v"
$ws.Range("D3").Value = 0.96
$ws.Range("E3").Value = 9.8
$ws.Range("F3").Value = 0.6

$ws.Range("A4").Value = "-4.49 + x + 0.5y"
$ws.Range("B4").Value = -3.51
$ws.Range("C4").Value = "J_0_LP_v"
$ws.Range("D4").Value = 0.88
$ws.Range("E4").Value = 8.4
$ws.Range("F4").Value = 0

$ws.Range("A5").Value = "-11.809999999999999 + x - 2y"
$ws.Range("B5").Value = -9.809999999999999
$ws.Range("C5").Value = "J_Ne_L0_v"
$ws.Range("D5").Value = 0.29
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 6.0

$ws.Range("A6").Value = "-4.92 - y"
$ws.Range("B6").Value = -4.92
$ws.Range("C6").Value = "J_Ne_L0_v"
$ws.Range("D6").Value = 0.34
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 6.6000000000000005

# --- Sheet: Punto_modificado ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("A2").Value = 2.0300000000000002
$ws.Range("B2").Value = 4.92

# --- Sheet: Vector_bf ---
$ws = $wb.Worksheets.Item(5)
$ws.Range("A2").Value = 0.8300000000000001

# --- Sheet: Vector_BF ---
$ws = $wb.Worksheets.Item(6)
$ws.Range("A2").Value = -6.550000000000001
$ws.Range("A3").Value = -14.8
